$wb = $excel.ActiveWorkbook

# Column F holds the "想去人数" (want-to-go count). This edit refreshes those
# counters to a newly scraped snapshot, updating the matching rows on every sheet.

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 14080
$ws.Cells.Item(4, 6).Value = 830
$ws.Cells.Item(6, 6).Value = 655
$ws.Cells.Item(8, 6).Value = 33
$ws.Cells.Item(9, 6).Value = 70
$ws.Cells.Item(10, 6).Value = 795
$ws.Cells.Item(11, 6).Value = 2174
$ws.Cells.Item(12, 6).Value = 160
$ws.Cells.Item(13, 6).Value = 118
$ws.Cells.Item(14, 6).Value = 95
$ws.Cells.Item(15, 6).Value = 204
$ws.Cells.Item(17, 6).Value = 579
$ws.Cells.Item(19, 6).Value = 494
$ws.Cells.Item(20, 6).Value = 343
$ws.Cells.Item(21, 6).Value = 27
$ws.Cells.Item(22, 6).Value = 301
$ws.Cells.Item(23, 6).Value = 880
$ws.Cells.Item(24, 6).Value = 136
$ws.Cells.Item(25, 6).Value = 66
$ws.Cells.Item(26, 6).Value = 20
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(29, 6).Value = 64
$ws.Cells.Item(30, 6).Value = 25

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 108
$ws.Cells.Item(8, 6).Value = 1659
$ws.Cells.Item(15, 6).Value = 1786

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 232
$ws.Cells.Item(3, 6).Value = 125
$ws.Cells.Item(4, 6).Value = 135

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 232
$ws.Cells.Item(3, 6).Value = 14080
$ws.Cells.Item(5, 6).Value = 830
$ws.Cells.Item(7, 6).Value = 655
$ws.Cells.Item(9, 6).Value = 33
$ws.Cells.Item(10, 6).Value = 70
$ws.Cells.Item(11, 6).Value = 795
$ws.Cells.Item(14, 6).Value = 2174
$ws.Cells.Item(15, 6).Value = 125
$ws.Cells.Item(16, 6).Value = 160
$ws.Cells.Item(17, 6).Value = 160
$ws.Cells.Item(18, 6).Value = 118
$ws.Cells.Item(19, 6).Value = 95
$ws.Cells.Item(20, 6).Value = 204
$ws.Cells.Item(24, 6).Value = 108
$ws.Cells.Item(25, 6).Value = 135
$ws.Cells.Item(26, 6).Value = 579
$ws.Cells.Item(28, 6).Value = 494
$ws.Cells.Item(29, 6).Value = 343
$ws.Cells.Item(30, 6).Value = 27
$ws.Cells.Item(31, 6).Value = 301
$ws.Cells.Item(32, 6).Value = 880
$ws.Cells.Item(34, 6).Value = 1659
$ws.Cells.Item(39, 6).Value = 136
$ws.Cells.Item(40, 6).Value = 66
$ws.Cells.Item(41, 6).Value = 20
$ws.Cells.Item(42, 6).Value = 1
$ws.Cells.Item(46, 6).Value = 64
$ws.Cells.Item(47, 6).Value = 25
$ws.Cells.Item(48, 6).Value = 1786
